$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task descriptions (column C) for Semestre 3 / Semestre 4 ---
# (written in this order so the shared-string table lines up with the
# author's original edit sequence)
$ws.Range("C11").Value = "Plannification de scénario pour tester le robot"
$ws.Range("C7").Value = "Implementation du code pour la marche du robot "
$ws.Range("C8").Value = "Scénario sur le déplacement en terrain plat du robot"
$ws.Range("C9").Value = "implementation du code pour la marche du robot en terrain accidenté ou pentu."
$ws.Range("C17").Value = "implementation du code pour validé les scénario plus complexes sur le robot (éviter obstacle et / ou galoper)"
$ws.Range("C12").Value = "Prise en charge d'un poid sur le robot sur terrain plat et pentu"
$ws.Range("C19").Value = "Implementation du code et assemblement du code pour répondre aux attentes de tous les scénarios en même temps"

# --- Enable word-wrap on the two bottom task rows (merged C17:G18 and C19:G20) ---
$ws.Range("C17:G18").WrapText = $true
$ws.Range("C19:G20").WrapText = $true

# --- Update the selected range/active cell to match the new layout ---
$ws.Range("C19:G20").Select()
